$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 144-146, pushing existing rows 144:211 down to 147:214
$ws.Range("A144:A146").EntireRow.Insert()

# Data common to all 3 newly inserted rows
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonD = 45029
$commonE = 13
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103002
$commonJ = "Ciruela"

# Row 144: Angeleno / Especial
$ws.Cells.Item(144, 1).Value = $commonA
$ws.Cells.Item(144, 2).Value = $commonB
$ws.Cells.Item(144, 3).Value = $commonC
$ws.Cells.Item(144, 4).Value = $commonD
$ws.Cells.Item(144, 5).Value = $commonE
$ws.Cells.Item(144, 6).Value = $commonF
$ws.Cells.Item(144, 7).Value = $commonG
$ws.Cells.Item(144, 8).Value = $commonH
$ws.Cells.Item(144, 9).Value = $commonI
$ws.Cells.Item(144, 10).Value = $commonJ
$ws.Cells.Item(144, 11).Value = "Angeleno"
$ws.Cells.Item(144, 12).Value = "Especial"
$ws.Cells.Item(144, 13).Value = 250
$ws.Cells.Item(144, 14).Value = 12000
$ws.Cells.Item(144, 15).Value = 12000
$ws.Cells.Item(144, 16).Value = 12000
$ws.Cells.Item(144, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(144, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(144, 19).Value = 800
$ws.Cells.Item(144, 20).Value = 15

# Row 145: Angeleno / Primera
$ws.Cells.Item(145, 1).Value = $commonA
$ws.Cells.Item(145, 2).Value = $commonB
$ws.Cells.Item(145, 3).Value = $commonC
$ws.Cells.Item(145, 4).Value = $commonD
$ws.Cells.Item(145, 5).Value = $commonE
$ws.Cells.Item(145, 6).Value = $commonF
$ws.Cells.Item(145, 7).Value = $commonG
$ws.Cells.Item(145, 8).Value = $commonH
$ws.Cells.Item(145, 9).Value = $commonI
$ws.Cells.Item(145, 10).Value = $commonJ
$ws.Cells.Item(145, 11).Value = "Angeleno"
$ws.Cells.Item(145, 12).Value = "Primera"
$ws.Cells.Item(145, 13).Value = 220
$ws.Cells.Item(145, 14).Value = 9000
$ws.Cells.Item(145, 15).Value = 9000
$ws.Cells.Item(145, 16).Value = 9000
$ws.Cells.Item(145, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(145, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(145, 19).Value = 600
$ws.Cells.Item(145, 20).Value = 15

# Row 146: Angeleno / Segunda
$ws.Cells.Item(146, 1).Value = $commonA
$ws.Cells.Item(146, 2).Value = $commonB
$ws.Cells.Item(146, 3).Value = $commonC
$ws.Cells.Item(146, 4).Value = $commonD
$ws.Cells.Item(146, 5).Value = $commonE
$ws.Cells.Item(146, 6).Value = $commonF
$ws.Cells.Item(146, 7).Value = $commonG
$ws.Cells.Item(146, 8).Value = $commonH
$ws.Cells.Item(146, 9).Value = $commonI
$ws.Cells.Item(146, 10).Value = $commonJ
$ws.Cells.Item(146, 11).Value = "Angeleno"
$ws.Cells.Item(146, 12).Value = "Segunda"
$ws.Cells.Item(146, 13).Value = 180
$ws.Cells.Item(146, 14).Value = 6000
$ws.Cells.Item(146, 15).Value = 6000
$ws.Cells.Item(146, 16).Value = 6000
$ws.Cells.Item(146, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(146, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(146, 19).Value = 400
$ws.Cells.Item(146, 20).Value = 15
